$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.6687250000000001
$ws.Range("H2").Value = 2.006175
$ws.Range("I2").Value = 0.3866332538806421
$ws.Range("J2").Value = 0.3866332538806422
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.07429999999999999
$ws.Range("N2").Value = 0.2229
$ws.Range("O2").Value = 0.08480128194885443
$ws.Range("P2").Value = 0.08480128194885443
$ws.Range("Q2").Value = 0.04968626750000001
$ws.Range("R2").Value = 0.4471764075
$ws.Range("S2").Value = 0.03278699557313535
$ws.Range("T2").Value = 0.03278699557313535

# Row 3
$ws.Range("G3").Value = 0.6687250000000001
$ws.Range("H3").Value = 2.006175
$ws.Range("I3").Value = 0.3866332538806421
$ws.Range("J3").Value = 0.3866332538806422
$ws.Range("O3").Value = 0.7166537695672586
$ws.Range("P3").Value = 0.7166537695672586
$ws.Range("Q3").Value = 0.4198975543916667
$ws.Range("R3").Value = 3.779077989525001
$ws.Range("S3").Value = 0.2770821788336171
$ws.Range("T3").Value = 0.2770821788336171

# Row 4
$ws.Range("G4").Value = 0.6687250000000001
$ws.Range("H4").Value = 2.006175
$ws.Range("I4").Value = 0.3866332538806421
$ws.Range("J4").Value = 0.3866332538806422
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.1739583333333333
$ws.Range("N4").Value = 0.521875
$ws.Range("O4").Value = 0.198544948483887
$ws.Range("P4").Value = 0.198544948483887
$ws.Range("Q4").Value = 0.1163302864583333
$ws.Range("R4").Value = 1.046972578125
$ws.Range("S4").Value = 0.07676407947388969
$ws.Range("T4").Value = 0.0767640794738897

# Row 5
$ws.Range("I5").Value = 0.5423686872113029
$ws.Range("J5").Value = 0.542368687211303
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.07429999999999999
$ws.Range("N5").Value = 0.2229
$ws.Range("O5").Value = 0.08480128194885443
$ws.Range("P5").Value = 0.08480128194885443
$ws.Range("Q5").Value = 0.06969983933333332
$ws.Range("R5").Value = 0.627298554
$ws.Range("S5").Value = 0.04599355996443574
$ws.Range("T5").Value = 0.04599355996443575

# Row 6
$ws.Range("I6").Value = 0.5423686872113029
$ws.Range("J6").Value = 0.542368687211303
$ws.Range("O6").Value = 0.7166537695672586
$ws.Range("P6").Value = 0.7166537695672586
$ws.Range("S6").Value = 0.3886905641852256
$ws.Range("T6").Value = 0.3886905641852257

# Row 7
$ws.Range("I7").Value = 0.5423686872113029
$ws.Range("J7").Value = 0.542368687211303
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.1739583333333333
$ws.Range("N7").Value = 0.521875
$ws.Range("O7").Value = 0.198544948483887
$ws.Range("P7").Value = 0.198544948483887
$ws.Range("Q7").Value = 0.1631879930555555
$ws.Range("R7").Value = 1.4686919375
$ws.Range("S7").Value = 0.1076845630616416
$ws.Range("T7").Value = 0.1076845630616416

# Row 8
$ws.Range("G8").Value = 0.122799
$ws.Range("H8").Value = 0.368397
$ws.Range("I8").Value = 0.07099805890805483
$ws.Range("J8").Value = 0.07099805890805484
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.07429999999999999
$ws.Range("N8").Value = 0.2229
$ws.Range("O8").Value = 0.08480128194885443
$ws.Range("P8").Value = 0.08480128194885443
$ws.Range("Q8").Value = 0.009123965699999998
$ws.Range("R8").Value = 0.08211569129999999
$ws.Range("S8").Value = 0.006020726411283333
$ws.Range("T8").Value = 0.006020726411283334

# Row 9
$ws.Range("G9").Value = 0.122799
$ws.Range("H9").Value = 0.368397
$ws.Range("I9").Value = 0.07099805890805483
$ws.Range("J9").Value = 0.07099805890805484
$ws.Range("O9").Value = 0.7166537695672586
$ws.Range("P9").Value = 0.7166537695672586
$ws.Range("Q9").Value = 0.07710643355899999
$ws.Range("R9").Value = 0.6939579020309999
$ws.Range("S9").Value = 0.05088102654841577
$ws.Range("T9").Value = 0.05088102654841579

# Row 10
$ws.Range("G10").Value = 0.122799
$ws.Range("H10").Value = 0.368397
$ws.Range("I10").Value = 0.07099805890805483
$ws.Range("J10").Value = 0.07099805890805484
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.1739583333333333
$ws.Range("N10").Value = 0.521875
$ws.Range("O10").Value = 0.198544948483887
$ws.Range("P10").Value = 0.198544948483887
$ws.Range("Q10").Value = 0.021361909375
$ws.Range("R10").Value = 0.192257184375
$ws.Range("S10").Value = 0.01409630594835572
$ws.Range("T10").Value = 0.01409630594835572
